$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 9 (Técnico Subsequente em Metalurgia - Campus Ouro Branco): Inscritos 10 -> 11
$ws.Range("E9").Value = 11

# Row 15 (Técnico Subsequente em Segurança do Trabalho - Campus Ouro Preto):
# Inscritos 54 -> 56, Pagos 25 -> 26, Inscrições homologadas 25 -> 26
$ws.Range("E15").Value = 56
$ws.Range("F15").Value = 26
$ws.Range("H15").Value = 26

# Row 16 (Técnico Subsequente em Logística - Campus Ribeirão das Neves): Inscritos 206 -> 207
$ws.Range("E16").Value = 207

# Row 18 (Técnico Subsequente em Segurança do Trabalho - Campus Santa Luzia): Inscritos 48 -> 51
$ws.Range("E18").Value = 51
